# Add a new slide ("Task:" slide) right after the existing slide, using
# the same master/layout ("Blank Slide" on the 2nd design) as slide 1,
# then populate it with two text boxes that mirror the geometry/style of
# the shapes on slide 1 ("CustomShape 1" title, "CustomShape 2" body).

$EMU_PER_PT = 12700.0

$p = $ppt.ActivePresentation

# Slide 1 uses the "Blank Slide" layout belonging to the 2nd slide master
# (design). Re-use the exact same layout for the new slide.
$design = $p.Designs.Item(2)
$layout = $design.SlideMaster.CustomLayouts.Item(1)
$slide2 = $p.Slides.AddSlide(2, $layout)

# ---- Shape 1: title textbox ("CustomShape 1") ----
$titleLeft = 497160 / $EMU_PER_PT
$titleTop = 57960 / $EMU_PER_PT
$titleWidth = 7820280 / $EMU_PER_PT
$titleHeight = 1160640 / $EMU_PER_PT

$shp1 = $slide2.Shapes.AddTextbox(1, $titleLeft, $titleTop, $titleWidth, $titleHeight)
$shp1.Name = "CustomShape 1"

$tf1 = $shp1.TextFrame
$tf1.MarginLeft = 90000 / $EMU_PER_PT
$tf1.MarginTop = 45000 / $EMU_PER_PT
$tf1.MarginRight = 90000 / $EMU_PER_PT
$tf1.MarginBottom = 45000 / $EMU_PER_PT
$tf1.VerticalAnchor = 4

$tr1 = $tf1.TextRange
$tr1.Text = "Default Methods"
$tr1.Font.Size = 32
$tr1.Font.Color.RGB = 0x927F5C
$tr1.Font.Name = "Lucida Bright"

$pr1 = $tr1.Paragraphs(1, 1)
$pr1.ParagraphFormat.SpaceWithin = 1

# ---- Shape 2: body textbox ("CustomShape 2") ----
$bodyLeft = 871200 / $EMU_PER_PT
$bodyTop = 1626840 / $EMU_PER_PT
$bodyWidth = 7808400 / $EMU_PER_PT
$bodyHeight = 4466456 / $EMU_PER_PT

$shp2 = $slide2.Shapes.AddTextbox(1, $bodyLeft, $bodyTop, $bodyWidth, $bodyHeight)
$shp2.Name = "CustomShape 2"

$tf2 = $shp2.TextFrame
$tf2.MarginLeft = 90000 / $EMU_PER_PT
$tf2.MarginTop = 45000 / $EMU_PER_PT
$tf2.MarginRight = 90000 / $EMU_PER_PT
$tf2.MarginBottom = 45000 / $EMU_PER_PT

$tr2 = $tf2.TextRange
$tr2.Text = "`rTask:`rCreate a new Interface with one abstract method and one default method`rCreate a new class implementing your interface`rCall both your implemented method and the default method"
$tr2.Font.Size = 20
$tr2.Font.Color.RGB = 0xB49530
$tr2.LanguageID = "sv-SE"

$count2 = $tr2.Paragraphs().Count

for ($i = 1; $i -le $count2; $i++) {
    $para = $tr2.Paragraphs($i, 1)
    $para.ParagraphFormat.SpaceWithin = 1
}

# Paragraphs 3, 4 and 5 ("Create a new Interface ...", "Create a new class
# ...", "Call both your implemented ...") are the numbered task list.
foreach ($i in 3, 4, 5) {
    $para = $tr2.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Type = 2
    $para.ParagraphFormat.Bullet.Style = 3
}
